$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above the current row 668, shifting rows 668:731
# down to 669:732 and extending the used range to A1:R732.
$ws.Rows.Item(668).Insert()

# Populate the newly inserted row 668 with the new weekly price record.
$ws.Cells.Item(668, 1).Value = 6
$ws.Cells.Item(668, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(668, 3).Value = "Metropolitana"
$ws.Cells.Item(668, 4).Value = 45132
$ws.Cells.Item(668, 5).Value = 13
$ws.Cells.Item(668, 6).Value = 100112039
$ws.Cells.Item(668, 7).Value = "Ciboulette"
$ws.Cells.Item(668, 8).Value = "Sin especificar"
$ws.Cells.Item(668, 9).Value = "Primera"
$ws.Cells.Item(668, 10).Value = 640
$ws.Cells.Item(668, 11).Value = 900
$ws.Cells.Item(668, 12).Value = 1000
$ws.Cells.Item(668, 13).Value = 956
$ws.Cells.Item(668, 14).Value = "`$/docena de atados"
$ws.Cells.Item(668, 15).Value = "Región Metropolitana"
$ws.Cells.Item(668, 16).Value = 319
$ws.Cells.Item(668, 17).Value = 3
$ws.Cells.Item(668, 18).Value = "Hortaliza"
